$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 20 de Marzo de 2020 a las 10:16'

$ws.Cells.Item(8, 2).Value = 15439
$ws.Cells.Item(8, 3).Value = 119
$ws.Cells.Item(8, 5).Value = 15280

$ws.Cells.Item(17, 2).Value = 1791
$ws.Cells.Item(17, 3).Value = 1
$ws.Cells.Item(17, 5).Value = 1783

$ws.Cells.Item(19, 2).Value = 1226
$ws.Cells.Item(19, 3).Value = 75
$ws.Cells.Item(19, 5).Value = 1216
$ws.Cells.Item(19, 7).Value = 3
$ws.Cells.Item(19, 8).Value = 9

$ws.Cells.Item(25, 2).Value = 774
$ws.Cells.Item(25, 3).Value = 80
$ws.Cells.Item(25, 5).Value = 771

$ws.Cells.Item(29, 1).Value = 'Irlanda'
$ws.Cells.Item(29, 2).Value = 557
$ws.Cells.Item(29, 3).Value = 0
$ws.Cells.Item(29, 4).Value = 5
$ws.Cells.Item(29, 5).Value = 549
$ws.Cells.Item(29, 6).Value = 6
$ws.Cells.Item(29, 8).Value = 3

$ws.Cells.Item(30, 1).Value = 'Luxemburgo'
$ws.Cells.Item(30, 2).Value = 484
$ws.Cells.Item(30, 3).Value = 149
$ws.Cells.Item(30, 4).Value = 6
$ws.Cells.Item(30, 5).Value = 474
$ws.Cells.Item(30, 6).Value = 1
$ws.Cells.Item(30, 8).Value = 4

$ws.Cells.Item(35, 1).Value = 'Indonesia'
$ws.Cells.Item(35, 2).Value = 369
$ws.Cells.Item(35, 3).Value = 60
$ws.Cells.Item(35, 4).Value = 17
$ws.Cells.Item(35, 5).Value = 320
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(35, 7).Value = 7
$ws.Cells.Item(35, 8).Value = 32

$ws.Cells.Item(36, 1).Value = 'Polonia'
$ws.Cells.Item(36, 2).Value = 367
$ws.Cells.Item(36, 3).Value = 12
$ws.Cells.Item(36, 4).Value = 13
$ws.Cells.Item(36, 5).Value = 349
$ws.Cells.Item(36, 6).Value = 3
$ws.Cells.Item(36, 8).Value = 5

$ws.Cells.Item(37, 1).Value = 'Turquia'
$ws.Cells.Item(37, 2).Value = 359
$ws.Cells.Item(37, 4).Value = 0
$ws.Cells.Item(37, 5).Value = 355
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(37, 8).Value = 4

$ws.Cells.Item(38, 1).Value = 'Singapur'
$ws.Cells.Item(38, 2).Value = 345
$ws.Cells.Item(38, 4).Value = 124
$ws.Cells.Item(38, 5).Value = 221
$ws.Cells.Item(38, 6).Value = 14

$ws.Cells.Item(39, 1).Value = 'Chile'
$ws.Cells.Item(39, 2).Value = 342
$ws.Cells.Item(39, 4).Value = 0
$ws.Cells.Item(39, 5).Value = 342
$ws.Cells.Item(39, 6).Value = 6

$ws.Cells.Item(40, 1).Value = 'Islandia'
$ws.Cells.Item(40, 2).Value = 330
$ws.Cells.Item(40, 3).Value = 0
$ws.Cells.Item(40, 4).Value = 5
$ws.Cells.Item(40, 5).Value = 325
$ws.Cells.Item(40, 8).Value = 0

$ws.Cells.Item(41, 1).Value = 'Tailandia'
$ws.Cells.Item(41, 2).Value = 322
$ws.Cells.Item(41, 3).Value = 50
$ws.Cells.Item(41, 4).Value = 42
$ws.Cells.Item(41, 5).Value = 279
$ws.Cells.Item(41, 6).Value = 1

$ws.Cells.Item(42, 1).Value = 'Eslovenia'
$ws.Cells.Item(42, 2).Value = 319
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 5).Value = 318
$ws.Cells.Item(42, 6).Value = 6
$ws.Cells.Item(42, 8).Value = 1

$ws.Cells.Item(49, 1).Value = 'Hong Kong'
$ws.Cells.Item(49, 2).Value = 256
$ws.Cells.Item(49, 3).Value = 48
$ws.Cells.Item(49, 4).Value = 98
$ws.Cells.Item(49, 5).Value = 154
$ws.Cells.Item(49, 6).Value = 4
$ws.Cells.Item(49, 7).Value = 0
$ws.Cells.Item(49, 8).Value = 4

$ws.Cells.Item(50, 1).Value = 'Peru'
$ws.Cells.Item(50, 2).Value = 234
$ws.Cells.Item(50, 3).Value = 0
$ws.Cells.Item(50, 4).Value = 1
$ws.Cells.Item(50, 5).Value = 230
$ws.Cells.Item(50, 6).Value = 7
$ws.Cells.Item(50, 7).Value = 2
$ws.Cells.Item(50, 8).Value = 3

$ws.Cells.Item(51, 1).Value = 'Filipinas'
$ws.Cells.Item(51, 2).Value = 230
$ws.Cells.Item(51, 3).Value = 13
$ws.Cells.Item(51, 4).Value = 8
$ws.Cells.Item(51, 5).Value = 204
$ws.Cells.Item(51, 6).Value = 1
$ws.Cells.Item(51, 7).Value = 1
$ws.Cells.Item(51, 8).Value = 18

$ws.Cells.Item(63, 5).Value = 105
$ws.Cells.Item(63, 7).Value = 1
$ws.Cells.Item(63, 8).Value = 2

$ws.Cells.Item(74, 1).Value = 'Vietnam'
$ws.Cells.Item(74, 2).Value = 87
$ws.Cells.Item(74, 3).Value = 2
$ws.Cells.Item(74, 4).Value = 16
$ws.Cells.Item(74, 5).Value = 71
$ws.Cells.Item(74, 6).Value = 0
$ws.Cells.Item(74, 8).Value = 0

$ws.Cells.Item(75, 1).Value = 'Hungria'
$ws.Cells.Item(75, 3).Value = 12
$ws.Cells.Item(75, 4).Value = 7
$ws.Cells.Item(75, 5).Value = 77
$ws.Cells.Item(75, 6).Value = 6
$ws.Cells.Item(75, 8).Value = 1

$ws.Cells.Item(83, 1).Value = 'Marruecos'
$ws.Cells.Item(83, 2).Value = 66
$ws.Cells.Item(83, 3).Value = 3
$ws.Cells.Item(83, 4).Value = 2
$ws.Cells.Item(83, 5).Value = 61
$ws.Cells.Item(83, 6).Value = 1
$ws.Cells.Item(83, 7).Value = 1
$ws.Cells.Item(83, 8).Value = 3

$ws.Cells.Item(85, 1).Value = 'Albania'
$ws.Cells.Item(85, 2).Value = 64
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 62
$ws.Cells.Item(85, 6).Value = 2

$ws.Cells.Item(88, 1).Value = 'Tunez'
$ws.Cells.Item(88, 2).Value = 50
$ws.Cells.Item(88, 3).Value = 11
$ws.Cells.Item(88, 4).Value = 1
$ws.Cells.Item(88, 5).Value = 48
$ws.Cells.Item(88, 6).Value = 2
$ws.Cells.Item(88, 8).Value = 1

$ws.Cells.Item(89, 1).Value = 'Kazajistan'
$ws.Cells.Item(89, 3).Value = 5
$ws.Cells.Item(89, 4).Value = 0
$ws.Cells.Item(89, 5).Value = 49
$ws.Cells.Item(89, 6).Value = 0
$ws.Cells.Item(89, 8).Value = 0

$ws.Cells.Item(90, 1).Value = 'Moldavia'
$ws.Cells.Item(90, 2).Value = 49
$ws.Cells.Item(90, 6).Value = 3
$ws.Cells.Item(90, 8).Value = 1

$ws.Cells.Item(91, 1).Value = 'Lituania'
$ws.Cells.Item(91, 4).Value = 1
$ws.Cells.Item(91, 5).Value = 47
$ws.Cells.Item(91, 6).Value = 1

$ws.Cells.Item(92, 1).Value = 'Oman'
$ws.Cells.Item(92, 2).Value = 48
$ws.Cells.Item(92, 4).Value = 13
$ws.Cells.Item(92, 5).Value = 35

$ws.Cells.Item(93, 1).Value = 'Estado de Palestina'
$ws.Cells.Item(93, 3).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 47

$ws.Cells.Item(94, 1).Value = 'Camboya'
$ws.Cells.Item(94, 2).Value = 47
$ws.Cells.Item(94, 3).Value = 10
$ws.Cells.Item(94, 4).Value = 1
$ws.Cells.Item(94, 5).Value = 46

$ws.Cells.Item(95, 1).Value = 'Guadalupe'
$ws.Cells.Item(95, 2).Value = 45
$ws.Cells.Item(95, 3).Value = 12
$ws.Cells.Item(95, 4).Value = 0
$ws.Cells.Item(95, 5).Value = 45
$ws.Cells.Item(95, 8).Value = 0

$ws.Cells.Item(96, 1).Value = 'Azerbaiyan'
$ws.Cells.Item(96, 2).Value = 44
$ws.Cells.Item(96, 3).Value = 0
$ws.Cells.Item(96, 4).Value = 7
$ws.Cells.Item(96, 5).Value = 36
$ws.Cells.Item(96, 6).Value = 0
$ws.Cells.Item(96, 8).Value = 1

$ws.Cells.Item(97, 1).Value = 'Georgia'
$ws.Cells.Item(97, 2).Value = 43
$ws.Cells.Item(97, 3).Value = 3
$ws.Cells.Item(97, 4).Value = 1
$ws.Cells.Item(97, 6).Value = 1

$ws.Cells.Item(98, 1).Value = 'Venezuela'
$ws.Cells.Item(98, 2).Value = 42
$ws.Cells.Item(98, 3).Value = 0
$ws.Cells.Item(98, 5).Value = 42

$ws.Cells.Item(99, 1).Value = 'Nueva Zelanda'
$ws.Cells.Item(99, 3).Value = 11
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 39
$ws.Cells.Item(99, 6).Value = 0
$ws.Cells.Item(99, 8).Value = 0

$ws.Cells.Item(105, 1).Value = 'Reunion'

$ws.Cells.Item(106, 1).Value = 'Liechtenstein'

$ws.Cells.Item(110, 1).Value = 'Camerun'
$ws.Cells.Item(110, 2).Value = 20
$ws.Cells.Item(110, 3).Value = 7
$ws.Cells.Item(110, 4).Value = 2

$ws.Cells.Item(111, 1).Value = 'Consejo Danes para los Refugiados'
$ws.Cells.Item(111, 3).Value = 4
$ws.Cells.Item(111, 4).Value = 0
$ws.Cells.Item(111, 5).Value = 18
$ws.Cells.Item(111, 8).Value = 0

$ws.Cells.Item(112, 1).Value = 'Banglades'
$ws.Cells.Item(112, 2).Value = 18
$ws.Cells.Item(112, 3).Value = 0
$ws.Cells.Item(112, 4).Value = 3
$ws.Cells.Item(112, 5).Value = 14
$ws.Cells.Item(112, 8).Value = 1

$ws.Cells.Item(113, 1).Value = 'Bolivia'
$ws.Cells.Item(113, 3).Value = 2
$ws.Cells.Item(113, 4).Value = 0
$ws.Cells.Item(113, 5).Value = 17

$ws.Cells.Item(114, 1).Value = 'Macao'
$ws.Cells.Item(114, 2).Value = 17
$ws.Cells.Item(114, 3).Value = 0
$ws.Cells.Item(114, 4).Value = 10
$ws.Cells.Item(114, 5).Value = 7
$ws.Cells.Item(114, 8).Value = 0

$ws.Cells.Item(115, 1).Value = 'Cuba'
$ws.Cells.Item(115, 3).Value = 5
$ws.Cells.Item(115, 4).Value = 0
$ws.Cells.Item(115, 5).Value = 15

$ws.Cells.Item(116, 1).Value = 'Jamaica'
$ws.Cells.Item(116, 2).Value = 16
$ws.Cells.Item(116, 3).Value = 1
$ws.Cells.Item(116, 4).Value = 2
$ws.Cells.Item(116, 5).Value = 13
$ws.Cells.Item(116, 8).Value = 1

$ws.Cells.Item(117, 1).Value = 'Guayana Francesa'
$ws.Cells.Item(117, 2).Value = 15
$ws.Cells.Item(117, 5).Value = 15

$ws.Cells.Item(118, 1).Value = 'Maldivas'
$ws.Cells.Item(118, 6).Value = 0

$ws.Cells.Item(119, 1).Value = 'Paraguay'
$ws.Cells.Item(119, 6).Value = 1

$ws.Cells.Item(120, 1).Value = 'Montenegro'
$ws.Cells.Item(120, 4).Value = 0
$ws.Cells.Item(120, 5).Value = 13

$ws.Cells.Item(123, 1).Value = 'Ghana'
$ws.Cells.Item(123, 3).Value = 0

$ws.Cells.Item(124, 1).Value = 'Monaco'
$ws.Cells.Item(124, 3).Value = 1

$ws.Cells.Item(125, 1).Value = 'Ruanda'

$ws.Cells.Item(126, 1).Value = 'Polinesia Francesa'
$ws.Cells.Item(126, 2).Value = 11
$ws.Cells.Item(126, 3).Value = 5
$ws.Cells.Item(126, 5).Value = 11

$ws.Cells.Item(135, 1).Value = 'Tanzania'

$ws.Cells.Item(136, 1).Value = 'Mongolia'

$ws.Cells.Item(137, 1).Value = 'Guinea Ecuatorial'
$ws.Cells.Item(137, 3).Value = 0

$ws.Cells.Item(138, 1).Value = 'Puerto Rico'

$ws.Cells.Item(139, 1).Value = 'Kirguistan'
$ws.Cells.Item(139, 3).Value = 3

$ws.Cells.Item(140, 1).Value = 'Seychelles'

$ws.Cells.Item(142, 1).Value = 'Guyana'
$ws.Cells.Item(142, 4).Value = 0
$ws.Cells.Item(142, 8).Value = 1

$ws.Cells.Item(143, 1).Value = 'Aruba'
$ws.Cells.Item(143, 4).Value = 1
$ws.Cells.Item(143, 8).Value = 0

$ws.Cells.Item(146, 1).Value = 'Congo'

$ws.Cells.Item(147, 1).Value = 'San Bartolome'

$ws.Cells.Item(148, 1).Value = 'Bahamas'

$ws.Cells.Item(149, 1).Value = 'San Martin (Parte Francesa)'

$ws.Cells.Item(150, 1).Value = 'Islas Virgenes de los Estados Unidos'

$ws.Cells.Item(151, 1).Value = 'Namibia'

$ws.Cells.Item(154, 1).Value = 'Bermudas'

$ws.Cells.Item(155, 1).Value = 'Groenlandia'

$ws.Cells.Item(156, 1).Value = 'Mauritania'

$ws.Cells.Item(157, 1).Value = 'Zambia'

$ws.Cells.Item(158, 1).Value = 'Butan'
$ws.Cells.Item(158, 3).Value = 1

$ws.Cells.Item(159, 1).Value = 'Republica de Africa Central'
$ws.Cells.Item(159, 3).Value = 1

$ws.Cells.Item(160, 1).Value = 'Nueva Caledonia'

$ws.Cells.Item(161, 1).Value = 'Benin'
$ws.Cells.Item(161, 3).Value = 0

$ws.Cells.Item(162, 1).Value = 'Santa Lucia'
$ws.Cells.Item(162, 3).Value = 0

$ws.Cells.Item(163, 1).Value = 'Liberia'
$ws.Cells.Item(163, 3).Value = 0

$ws.Cells.Item(164, 1).Value = 'Haiti'
$ws.Cells.Item(164, 3).Value = 2

$ws.Cells.Item(166, 1).Value = 'Montserrat'

$ws.Cells.Item(167, 1).Value = 'Gambia'
$ws.Cells.Item(167, 3).Value = 0

$ws.Cells.Item(168, 1).Value = 'Isla de Man'

$ws.Cells.Item(169, 1).Value = 'Antigua y Barbuda'

$ws.Cells.Item(170, 1).Value = 'Cabo Verde'
$ws.Cells.Item(170, 3).Value = 1

$ws.Cells.Item(171, 1).Value = 'Surinam'

$ws.Cells.Item(172, 1).Value = 'Fiyi'

$ws.Cells.Item(173, 1).Value = 'Santa Sede'

$ws.Cells.Item(174, 1).Value = 'Guinea'

$ws.Cells.Item(175, 1).Value = 'Angola'
$ws.Cells.Item(175, 3).Value = 1

$ws.Cells.Item(176, 1).Value = 'Suazilandia'

$ws.Cells.Item(177, 1).Value = 'Somalia'

$ws.Cells.Item(178, 1).Value = 'San Martin (Parte Holandesa)'

$ws.Cells.Item(179, 1).Value = 'El Salvador'

$ws.Cells.Item(180, 1).Value = 'Nicaragua'

$ws.Cells.Item(181, 1).Value = 'Republica del Chad'

$ws.Cells.Item(182, 1).Value = 'Niger'
$ws.Cells.Item(182, 3).Value = 0

$ws.Cells.Item(183, 1).Value = 'Republica de Yibuti'

$ws.Cells.Item(184, 1).Value = 'San Vicente y las Granadinas'
